# Primer caso de prueba creado, sin errores.
# - Elimina la hoja "Hoja1" (ya no se necesita, solo queda "DatosCP").
# - Actualiza el primer caso de prueba (fila 2) en "DatosCP" con datos
#   nuevos relacionados al login con email invalido.

$wb = $excel.ActiveWorkbook

# Quita la hoja sobrante "Hoja1"; solo debe quedar "DatosCP".
$wb.Worksheets("Hoja1").Delete()

$ws = $wb.Worksheets("DatosCP")

# Nuevo primer caso de prueba (fila 2).
$ws.Range("B2").Value = "qweqweew"
$ws.Range("C2").Value = "ee51165"
$ws.Range("A2").Value = "CP001_loginInvalidEmail"
$ws.Range("D2").Value = "Invalid email address."

# Resalta el identificador del caso de prueba con subrayado.
$ws.Range("A2").Font.Underline = 2
